$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 249, shifting existing rows 249:262 down to 250:263
# (dimension grows from A1:R262 to A1:R263 automatically).
$ws.Rows.Item(249).Insert()

# Populate the newly inserted row 249 with the new record.
$ws.Cells.Item(249, 1).Value = 4
$ws.Cells.Item(249, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(249, 3).Value = "Los Lagos"
$ws.Cells.Item(249, 4).Value = 44610
$ws.Cells.Item(249, 5).Value = 10
$ws.Cells.Item(249, 6).Value = 100112045
$ws.Cells.Item(249, 7).Value = "Zapallo"
$ws.Cells.Item(249, 8).Value = "Paine"
$ws.Cells.Item(249, 9).Value = "1a nueva(o)"
$ws.Cells.Item(249, 10).Value = 1200
$ws.Cells.Item(249, 11).Value = 500
$ws.Cells.Item(249, 12).Value = 500
$ws.Cells.Item(249, 13).Value = 500
$ws.Cells.Item(249, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(249, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(249, 16).Value = 500
$ws.Cells.Item(249, 17).Value = 1
$ws.Cells.Item(249, 18).Value = "Hortaliza"
